# Scheduled runner update: refresh market-price / profit figures (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 786.1539
$ws.Range("I18").Value = 810.0833
$ws.Range("J18").Value = 499
$ws.Range("K18").Value = 810.0833
$ws.Range("L18").Value = 499
$ws.Range("M18").Value = -526.0833
$ws.Range("N18").Value = -1067
$ws.Range("H30").Value = 9000
$ws.Range("J30").Value = 9000
$ws.Range("L30").Value = 27000
$ws.Range("N30").Value = -27202
$ws.Range("H31").Value = 1681.6666
$ws.Range("J31").Value = 23.5
$ws.Range("L31").Value = 70.5
$ws.Range("N31").Value = -530.5
$ws.Range("H38").Value = 902.55554
$ws.Range("I38").Value = 1009.125
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 3027.375
$ws.Range("L38").Value = 150
$ws.Range("M38").Value = -2655.375
$ws.Range("N38").Value = -894
$ws.Range("H39").Value = 401.36667
$ws.Range("I39").Value = 728
$ws.Range("J39").Value = 282.5909
$ws.Range("K39").Value = 2184
$ws.Range("L39").Value = 847.7727
$ws.Range("M39").Value = -1888
$ws.Range("N39").Value = -1439.7727
$ws.Range("H42").Value = 232.76923
$ws.Range("I42").Value = 220.54546
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 661.6363799999999
$ws.Range("L42").Value = 900
$ws.Range("M42").Value = -431.6363799999999
$ws.Range("N42").Value = -1360
$ws.Range("H51").Value = 105151.1
$ws.Range("J51").Value = 8167.3335
$ws.Range("L51").Value = 8167.3335
$ws.Range("N51").Value = -9135.333500000001
$ws.Range("H59").Value = 1553.1666
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 1663.8
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 4991.4
$ws.Range("M59").Value = -2443
$ws.Range("N59").Value = -6105.4
$ws.Range("H112").Value = 1728.125
$ws.Range("J112").Value = 1963.6111
$ws.Range("L112").Value = 5890.8333
$ws.Range("N112").Value = -8106.8333
$ws.Range("H135").Value = 1273.5128
$ws.Range("I135").Value = 1394.8857
$ws.Range("J135").Value = 211.5
$ws.Range("K135").Value = 12553.9713
$ws.Range("L135").Value = 1903.5
$ws.Range("M135").Value = -10018.9713
$ws.Range("N135").Value = -6973.5
$ws.Range("H137").Value = 3418.1428
$ws.Range("I137").Value = 2450.3462
$ws.Range("K137").Value = 7351.0386
$ws.Range("M137").Value = -4801.0386
$ws.Range("H138").Value = 2454.7805
$ws.Range("J138").Value = 2761.5
$ws.Range("L138").Value = 8284.5
$ws.Range("N138").Value = -18564.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7261133
$ws.Range("I32").Value = 7706310.5
$ws.Range("J32").Value = 27003
$ws.Range("K32").Value = 7706310.5
$ws.Range("L32").Value = 27003
$ws.Range("M32").Value = -7706023.5
$ws.Range("N32").Value = -27577
$ws.Range("H132").Value = 2571.0527
$ws.Range("I132").Value = 1249.2273
$ws.Range("K132").Value = 3747.6819
$ws.Range("M132").Value = -1217.6819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 39565.332
$ws.Range("J76").Value = 41480.4
$ws.Range("L76").Value = 41480.4
$ws.Range("N76").Value = -42110.4
$ws.Range("H79").Value = 39565.332
$ws.Range("J79").Value = 41480.4
$ws.Range("L79").Value = 41480.4
$ws.Range("N79").Value = -43664.4
$ws.Range("H134").Value = 313541.56
$ws.Range("I134").Value = 924.85187
$ws.Range("K134").Value = 2774.55561
$ws.Range("M134").Value = -239.5556099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2168.6667
$ws.Range("I99").Value = 2045.5555
$ws.Range("J99").Value = 2538
$ws.Range("K99").Value = 2045.5555
$ws.Range("L99").Value = 2538
$ws.Range("M99").Value = -547.5554999999999
$ws.Range("N99").Value = -5534
$ws.Range("H105").Value = 1710
$ws.Range("I105").Value = 1710
$ws.Range("K105").Value = 1710
$ws.Range("M105").Value = 37
$ws.Range("H126").Value = 2168.6667
$ws.Range("I126").Value = 2045.5555
$ws.Range("J126").Value = 2538
$ws.Range("K126").Value = 6136.666499999999
$ws.Range("L126").Value = 7614
$ws.Range("M126").Value = -3666.666499999999
$ws.Range("N126").Value = -12554
$ws.Range("H132").Value = 1616.4082
$ws.Range("I132").Value = 1666.279
$ws.Range("J132").Value = 1259
$ws.Range("K132").Value = 4998.837
$ws.Range("L132").Value = 3777
$ws.Range("M132").Value = -2468.837
$ws.Range("N132").Value = -8837
$ws.Range("H134").Value = 2015.8636
$ws.Range("I134").Value = 1019.6667
$ws.Range("K134").Value = 3059.0001
$ws.Range("M134").Value = -524.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 3006.6667
$ws.Range("I13").Value = 1759.5
$ws.Range("J13").Value = 5501
$ws.Range("K13").Value = 5278.5
$ws.Range("L13").Value = 16503
$ws.Range("M13").Value = -5110.5
$ws.Range("N13").Value = -16839
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H37").Value = 85000
$ws.Range("J37").Value = 85000
$ws.Range("L37").Value = 255000
$ws.Range("N37").Value = -255224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1045.1428
$ws.Range("I22").Value = 1423.4
$ws.Range("J22").Value = 99.5
$ws.Range("K22").Value = 1423.4
$ws.Range("L22").Value = 99.5
$ws.Range("M22").Value = -894.4000000000001
$ws.Range("N22").Value = -1157.5
$ws.Range("H25").Value = 1005
$ws.Range("J25").Value = 1206.6666
$ws.Range("L25").Value = 1206.6666
$ws.Range("N25").Value = -2264.6666
$ws.Range("H63").Value = 21282.8
$ws.Range("J63").Value = 24103.5
$ws.Range("L63").Value = 24103.5
$ws.Range("N63").Value = -25475.5
$ws.Range("H66").Value = 21282.8
$ws.Range("J66").Value = 24103.5
$ws.Range("L66").Value = 72310.5
$ws.Range("N66").Value = -79174.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H7").Value = 45518.84
$ws.Range("J7").Value = 132338.88
$ws.Range("L7").Value = 132338.88
$ws.Range("N7").Value = -132562.88
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H40").Value = 2097.4167
$ws.Range("I40").Value = 1343.3793
$ws.Range("K40").Value = 1343.3793
$ws.Range("M40").Value = -1207.3793
$ws.Range("H56").Value = 25498.334
$ws.Range("I56").Value = 16000
$ws.Range("K56").Value = 16000
$ws.Range("M56").Value = -15309
$ws.Range("H122").Value = 5308.5
$ws.Range("I122").Value = 4908.1
$ws.Range("J122").Value = 6109.3
$ws.Range("K122").Value = 14724.3
$ws.Range("L122").Value = 18327.9
$ws.Range("M122").Value = -12274.3
$ws.Range("N122").Value = -23227.9
$ws.Range("H126").Value = 45518.84
$ws.Range("J126").Value = 132338.88
$ws.Range("L126").Value = 397016.64
$ws.Range("N126").Value = -401956.64
$ws.Range("H136").Value = 63982.668
$ws.Range("I136").Value = 8683.714
$ws.Range("J136").Value = 174580.58
$ws.Range("K136").Value = 26051.142
$ws.Range("L136").Value = 523741.74
$ws.Range("M136").Value = -23501.142
$ws.Range("N136").Value = -528841.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 44999
$ws.Range("J55").Value = 44999
$ws.Range("L55").Value = 44999
$ws.Range("N55").Value = -45553
$ws.Range("H122").Value = 1863.9231
$ws.Range("I122").Value = 1863.9231
$ws.Range("K122").Value = 5591.7693
$ws.Range("M122").Value = -3141.7693
$ws.Range("H136").Value = 2868.077
$ws.Range("I136").Value = 2355.6562
$ws.Range("K136").Value = 7066.9686
$ws.Range("M136").Value = -4516.9686
